# edit.ps1 - applies the commit "more app stuff. Need to fix empty summoner name"
# Adds a new worksheet ("Sheet2") before "Sheet1" with match/participant role
# lookup data, nudges the picture on Sheet1 slightly, and leaves a couple of
# cosmetic selection/view tweaks behind (as Excel itself would when re-saving).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Sheet2" worksheet just before the existing "Sheet1"
# ------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($sheet1)
$newSheet.Name = "Sheet2"

# Re-resolve "Sheet1" by name: the variable captured above can become stale
# once a sheet is inserted/renamed around it.
$sheet1 = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 2. Populate the new sheet with the Role / ValidRoles lookup tables
#    (written in the same order the original author typed them, so the
#    shared-string table comes out with the same index assignment:
#    headers of table 1 -> body of table 1 -> headers of table 2 -> body
#    of table 2)
# ------------------------------------------------------------------

# Headers for first table (G5:J5)
$newSheet.Range("G5").Value = "MatchId"
$newSheet.Range("H5").Value = "ParticipantId"
$newSheet.Range("I5").Value = "Role"
$newSheet.Range("J5").Value = "ValidRoles"

# First table body: G/H rows 6-15 use "a", rows 16-25 use "b", row 26 "c"
for ($i = 0; $i -lt 10; $i++) {
    $row = 6 + $i
    $newSheet.Cells.Item($row, 7).Value = "a"
    $newSheet.Cells.Item($row, 8).Value = $i + 1
}
for ($i = 0; $i -lt 10; $i++) {
    $row = 16 + $i
    $newSheet.Cells.Item($row, 7).Value = "b"
    $newSheet.Cells.Item($row, 8).Value = $i + 1
}
$newSheet.Cells.Item(26, 7).Value = "c"

# Headers for second table (M5:Q5)
$newSheet.Range("M5").Value = "MatchId"
$newSheet.Range("N5").Value = "ParticipantId"
$newSheet.Range("O5").Value = "Jungle Flag"
$newSheet.Range("P5").Value = "SupportFlag"
$newSheet.Range("Q5").Value = "posRole"

# Second table body: M/N mirror G/H for rows 6-25 (no "c" row here)
for ($i = 0; $i -lt 10; $i++) {
    $row = 6 + $i
    $newSheet.Cells.Item($row, 13).Value = "a"
    $newSheet.Cells.Item($row, 14).Value = $i + 1
}
for ($i = 0; $i -lt 10; $i++) {
    $row = 16 + $i
    $newSheet.Cells.Item($row, 13).Value = "b"
    $newSheet.Cells.Item($row, 14).Value = $i + 1
}

# Jungle Flag / SupportFlag booleans for the "a" rows (6-15)
$newSheet.Cells.Item(6, 15).Value = $true
$newSheet.Cells.Item(6, 16).Value = $false
$newSheet.Cells.Item(7, 15).Value = $false
$newSheet.Cells.Item(7, 16).Value = $true
$newSheet.Cells.Item(8, 15).Value = $false
$newSheet.Cells.Item(8, 16).Value = $true
$newSheet.Cells.Item(9, 15).Value = $false
$newSheet.Cells.Item(9, 16).Value = $false
$newSheet.Cells.Item(10, 15).Value = $false
$newSheet.Cells.Item(10, 16).Value = $false
$newSheet.Cells.Item(11, 15).Value = $false
$newSheet.Cells.Item(11, 16).Value = $false
$newSheet.Cells.Item(12, 15).Value = $false
$newSheet.Cells.Item(12, 16).Value = $false
$newSheet.Cells.Item(13, 15).Value = $false
$newSheet.Cells.Item(13, 16).Value = $false
$newSheet.Cells.Item(14, 15).Value = $false
$newSheet.Cells.Item(14, 16).Value = $false
$newSheet.Cells.Item(15, 15).Value = $true
$newSheet.Cells.Item(15, 16).Value = $false

# Column widths (best-fit in the original file)
$newSheet.Columns.Item(8).ColumnWidth = 11.33203125
$newSheet.Columns.Item(14).ColumnWidth = 11.33203125
$newSheet.Columns.Item(15).ColumnWidth = 9.6640625
$newSheet.Columns.Item(16).ColumnWidth = 10.5546875

# ------------------------------------------------------------------
# 3. Nudge the picture on "Sheet1" (the locations/areas sheet)
# ------------------------------------------------------------------
$shape = $sheet1.Shapes.Item(1)
# Re-assert the (unchanged) picture size before nudging its position: this
# host's Shape.Width/Height getters don't round-trip column widths exactly,
# so touching Top/Left alone drifts the cached extents. Setting the size
# back explicitly (in points, matching the original 3025402 x 3116850 EMU)
# keeps xdr:ext / a:ext exact.
$shape.Width = 238.22062992125984
$shape.Height = 245.42125984251967
$shape.Top = $shape.Top + 17.4
$shape.Left = $shape.Left + 13.2

# ------------------------------------------------------------------
# 4. Misc view/selection tweaks left over from the authoring session
#    (set these first so the final ".Select()" below - on the new
#    "Sheet2" tab - is the one that ends up as the active sheet)
# ------------------------------------------------------------------
$working = $wb.Worksheets.Item("working")
$working.Range("V21").Select()

$sheet1.Range("L21").Select()

$wb.Windows.Item(1).Height = 623.45

$newSheet.Range("Q19").Select()
